$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reassign the header row values for the 12/16/19 timepoint collection sheet:
# A1: short_cage_id -> mouse_id
# B1: n_mice -> short_cage_id
# C1: total_cage_weight_(g) stays the same
$ws.Range("A1").Value = "mouse_id"
$ws.Range("B1").Value = "short_cage_id"
$ws.Range("C1").Value = "total_cage_weight_(g)"
